$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 281.41666
$ws.Range("I2").Value = 201.25
$ws.Range("J2").Value = 441.75
$ws.Range("K2").Value = 201.25
$ws.Range("L2").Value = 441.75
$ws.Range("M2").Value = -88.25
$ws.Range("N2").Value = -667.75
# Row 6
$ws.Range("H6").Value = 1356
$ws.Range("I6").Value = 250.25
$ws.Range("J6").Value = 2830.3333
$ws.Range("K6").Value = 750.75
$ws.Range("L6").Value = 8490.999899999999
$ws.Range("M6").Value = -638.75
$ws.Range("N6").Value = -8714.999899999999
# Row 51
$ws.Range("H51").Value = 4803.24
$ws.Range("I51").Value = 2300
$ws.Range("J51").Value = 5429.05
$ws.Range("K51").Value = 2300
$ws.Range("L51").Value = 5429.05
$ws.Range("M51").Value = -1816
$ws.Range("N51").Value = -6397.05
# Row 86
$ws.Range("H86").Value = 21036.45
$ws.Range("I86").Value = 632.8125
$ws.Range("J86").Value = 102651
$ws.Range("K86").Value = 632.8125
$ws.Range("L86").Value = 102651
$ws.Range("M86").Value = 490.1875
$ws.Range("N86").Value = -104897
# Row 89
$ws.Range("H89").Value = 21036.45
$ws.Range("I89").Value = 632.8125
$ws.Range("J89").Value = 102651
$ws.Range("K89").Value = 3164.0625
$ws.Range("L89").Value = 513255
$ws.Range("M89").Value = 2451.9375
$ws.Range("N89").Value = -524487
# Row 105
$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("N105").Value = 0
$ws.Range("L105").ClearContents()
# Row 129
$ws.Range("H129").Value = 23810360
$ws.Range("J129").Value = 917.56757
$ws.Range("L129").Value = 2752.70271
$ws.Range("N129").Value = -12752.70271

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1536
$ws.Range("I2").Value = 1197.1
$ws.Range("J2").Value = 2100.8333
$ws.Range("K2").Value = 1197.1
$ws.Range("L2").Value = 2100.8333
$ws.Range("M2").Value = -1084.1
$ws.Range("N2").Value = -2326.8333
# Row 8
$ws.Range("H8").Value = 10000000
$ws.Range("I8").Value = 10000000
$ws.Range("K8").Value = 10000000
$ws.Range("M8").Value = -9999856
# Row 45
$ws.Range("H45").Value = 4087.4285
$ws.Range("I45").Value = 1778
$ws.Range("J45").Value = 7166.6665
$ws.Range("K45").Value = 1778
$ws.Range("L45").Value = 7166.6665
$ws.Range("M45").Value = -1401
$ws.Range("N45").Value = -7920.6665
# Row 110
$ws.Range("H110").Value = 46332.09
$ws.Range("I110").Value = 72124.57000000001
$ws.Range("J110").Value = 1195.25
$ws.Range("K110").Value = 72124.57000000001
$ws.Range("L110").Value = 1195.25
$ws.Range("M110").Value = -70079.57000000001
$ws.Range("N110").Value = -5285.25
# Row 116
$ws.Range("H116").Value = 1536
$ws.Range("I116").Value = 1197.1
$ws.Range("J116").Value = 2100.8333
$ws.Range("K116").Value = 1197.1
$ws.Range("L116").Value = 2100.8333
$ws.Range("M116").Value = 1096.9
$ws.Range("N116").Value = -6688.8333

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1536
$ws.Range("I3").Value = 1197.1
$ws.Range("J3").Value = 2100.8333
$ws.Range("K3").Value = 1197.1
$ws.Range("L3").Value = 2100.8333
$ws.Range("M3").Value = -1083.1
$ws.Range("N3").Value = -2328.8333
# Row 11
$ws.Range("H11").Value = 3003.3333
$ws.Range("J11").Value = 3005
$ws.Range("L11").Value = 3005
$ws.Range("N11").Value = -3285
# Row 105
$ws.Range("H105").Value = 3828.5715
$ws.Range("I105").Value = 3760
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 3760
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -2013
$ws.Range("N105").Value = -7494
# Row 134
$ws.Range("H134").Value = 2543.2812
$ws.Range("I134").Value = 1793.1904
$ws.Range("J134").Value = 3975.2727
$ws.Range("K134").Value = 5379.5712
$ws.Range("L134").Value = 11925.8181
$ws.Range("M134").Value = -2844.5712
$ws.Range("N134").Value = -16995.8181

$ws = $wb.Worksheets.Item("CRP")
# Row 16
$ws.Range("H16").Value = 1716.55
$ws.Range("I16").Value = 1727
$ws.Range("J16").Value = 1692.1666
$ws.Range("K16").Value = 1727
$ws.Range("L16").Value = 1692.1666
$ws.Range("M16").Value = -1440
$ws.Range("N16").Value = -2266.1666
# Row 22
$ws.Range("H22").Value = 1274.6364
$ws.Range("I22").Value = 1274.6364
$ws.Range("K22").Value = 1274.6364
$ws.Range("M22").Value = -924.6364000000001
# Row 31
$ws.Range("H31").Value = 11559.344
$ws.Range("I31").Value = 1120
$ws.Range("J31").Value = 42877.375
$ws.Range("K31").Value = 1120
$ws.Range("L31").Value = 42877.375
$ws.Range("M31").Value = -825
$ws.Range("N31").Value = -43467.375
# Row 34
$ws.Range("H34").Value = 11559.344
$ws.Range("I34").Value = 1120
$ws.Range("J34").Value = 42877.375
$ws.Range("K34").Value = 1120
$ws.Range("L34").Value = 42877.375
$ws.Range("M34").Value = -918
$ws.Range("N34").Value = -43281.375
# Row 99
$ws.Range("H99").Value = 2919763.5
$ws.Range("I99").Value = 3557266.5
$ws.Range("J99").Value = 51000
$ws.Range("K99").Value = 3557266.5
$ws.Range("L99").Value = 51000
$ws.Range("M99").Value = -3555768.5
$ws.Range("N99").Value = -53996
# Row 105
$ws.Range("H105").Value = 3102.5
$ws.Range("I105").Value = 4003.3333
$ws.Range("K105").Value = 4003.3333
$ws.Range("M105").Value = -2256.3333
# Row 113
$ws.Range("H113").Value = 1716.55
$ws.Range("I113").Value = 1727
$ws.Range("J113").Value = 1692.1666
$ws.Range("K113").Value = 1727
$ws.Range("L113").Value = 1692.1666
$ws.Range("M113").Value = 443
$ws.Range("N113").Value = -6032.1666
# Row 126
$ws.Range("H126").Value = 2919763.5
$ws.Range("I126").Value = 3557266.5
$ws.Range("J126").Value = 51000
$ws.Range("K126").Value = 10671799.5
$ws.Range("L126").Value = 153000
$ws.Range("M126").Value = -10669329.5
$ws.Range("N126").Value = -157940

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1244.8214
$ws.Range("J5").Value = 2609.9
$ws.Range("L5").Value = 7829.700000000001
$ws.Range("N5").Value = -8053.700000000001
# Row 37
$ws.Range("H37").Value = 79992.25
$ws.Range("J37").Value = 79992.25
$ws.Range("L37").Value = 239976.75
$ws.Range("N37").Value = -240200.75
# Row 80
$ws.Range("H80").Value = 4716.222
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 4805.75
$ws.Range("K80").Value = 12000
$ws.Range("L80").Value = 14417.25
$ws.Range("M80").Value = -11064
$ws.Range("N80").Value = -16289.25
# Row 83
$ws.Range("H83").Value = 4716.222
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 4805.75
$ws.Range("K83").Value = 36000
$ws.Range("L83").Value = 43251.75
$ws.Range("M83").Value = -31320
$ws.Range("N83").Value = -52611.75
# Row 92
$ws.Range("H92").Value = 200214.6
$ws.Range("I92").Value = 333533.34
$ws.Range("J92").Value = 236.5
$ws.Range("K92").Value = 1000600.02
$ws.Range("L92").Value = 709.5
$ws.Range("M92").Value = -999352.02
$ws.Range("N92").Value = -3205.5
# Row 114
$ws.Range("H114").Value = 532.2143
$ws.Range("I114").Value = 364.5
$ws.Range("J114").Value = 599.3
$ws.Range("K114").Value = 1093.5
$ws.Range("L114").Value = 1797.9
$ws.Range("M114").Value = 2160.5
$ws.Range("N114").Value = -8305.9
# Row 131
$ws.Range("H131").Value = 1225.8772
$ws.Range("I131").Value = 403.125
$ws.Range("J131").Value = 1360.2041
$ws.Range("K131").Value = 1209.375
$ws.Range("L131").Value = 4080.6123
$ws.Range("M131").Value = 3830.625
$ws.Range("N131").Value = -14160.6123
# Row 133
$ws.Range("H133").Value = 2150.625
$ws.Range("I133").Value = 2344.2856
$ws.Range("K133").Value = 7032.8568
$ws.Range("M133").Value = -1972.8568
# Row 135
$ws.Range("H135").Value = 1244.8214
$ws.Range("J135").Value = 2609.9
$ws.Range("L135").Value = 23489.1
$ws.Range("N135").Value = -28559.1

$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 12901
$ws.Range("J5").Value = 12901
$ws.Range("L5").Value = 12901
$ws.Range("N5").Value = -13125
# Row 113
$ws.Range("H113").Value = 1637
$ws.Range("I113").Value = 1576.6666
$ws.Range("J113").Value = 1999
$ws.Range("K113").Value = 1576.6666
$ws.Range("L113").Value = 1999
$ws.Range("M113").Value = 593.3334
$ws.Range("N113").Value = -6339
# Row 126
$ws.Range("H126").Value = 2035.0741
$ws.Range("I126").Value = 1646.8572
$ws.Range("K126").Value = 4940.571599999999
$ws.Range("M126").Value = -2470.571599999999
# Row 140
$ws.Range("H140").Value = 54992.5
$ws.Range("J140").Value = 54992.5
$ws.Range("L140").Value = 54992.5
$ws.Range("N140").Value = -65352.5

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 2169.8572
$ws.Range("I61").Value = 1646
$ws.Range("J61").Value = 2868.3333
$ws.Range("K61").Value = 1646
$ws.Range("L61").Value = 2868.3333
$ws.Range("M61").Value = -1444
$ws.Range("N61").Value = -3272.3333
# Row 113
$ws.Range("H113").Value = 2169.8572
$ws.Range("I113").Value = 1646
$ws.Range("J113").Value = 2868.3333
$ws.Range("K113").Value = 1646
$ws.Range("L113").Value = 2868.3333
$ws.Range("M113").Value = 524
$ws.Range("N113").Value = -7208.3333

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 38707.668
$ws.Range("I122").Value = 43387.793
$ws.Range("J122").Value = 1266.6666
$ws.Range("K122").Value = 130163.379
$ws.Range("L122").Value = 3799.9998
$ws.Range("M122").Value = -127713.379
$ws.Range("N122").Value = -8699.9998
# Row 126
$ws.Range("H126").Value = 40743.96
$ws.Range("I126").Value = 60315.293
$ws.Range("J126").Value = 3775.889
$ws.Range("K126").Value = 180945.879
$ws.Range("L126").Value = 11327.667
$ws.Range("M126").Value = -178475.879
$ws.Range("N126").Value = -16267.667
# Row 136
$ws.Range("H136").Value = 4014.6482
$ws.Range("I136").Value = 627.7059
$ws.Range("K136").Value = 1883.1177
$ws.Range("M136").Value = 666.8822999999998
